$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1844109867082959
$ws.Cells.Item(2, 3).Value = 0.2434303829804775
$ws.Cells.Item(2, 4).Value = 0.7550547198149107
$ws.Cells.Item(2, 5).Value = 0.6915081552015099
$ws.Cells.Item(2, 6).Value = 0.4175641082963374
$ws.Cells.Item(2, 7).Value = 0.5230975067184298
$ws.Cells.Item(2, 8).Value = 0.6743329343542892
$ws.Cells.Item(2, 9).Value = 0.188602888553233
$ws.Cells.Item(2, 10).Value = 0.4647918628642795
$ws.Cells.Item(2, 11).Value = 0.2901314376394514

$ws.Cells.Item(3, 2).Value = 0.2658583626134861
$ws.Cells.Item(3, 3).Value = 0.7576582593773593
$ws.Cells.Item(3, 4).Value = 0.6694404695474946
$ws.Cells.Item(3, 5).Value = 0.4115821279919825
$ws.Cells.Item(3, 6).Value = 0.5177542565360187
$ws.Cells.Item(3, 7).Value = 0.663964400983923
$ws.Cells.Item(3, 8).Value = 0.1801255664296092
$ws.Cells.Item(3, 9).Value = 0.4570195545574687
$ws.Cells.Item(3, 10).Value = 0.2815074641064828
$ws.Cells.Item(3, 11).Value = 0.589805676679554

$ws.Cells.Item(4, 2).Value = 0.7790828308345794
$ws.Cells.Item(4, 3).Value = 0.7896183183467069
$ws.Cells.Item(4, 4).Value = 0.3223877835781994
$ws.Cells.Item(4, 5).Value = 0.4880272526874408
$ws.Cells.Item(4, 6).Value = 0.6725905210138331
$ws.Cells.Item(4, 7).Value = 0.154185393004174
$ws.Cells.Item(4, 8).Value = 0.4342409915292412
$ws.Cells.Item(4, 9).Value = 0.267889428957865
$ws.Cells.Item(4, 10).Value = 0.5712848794943781
$ws.Cells.Item(4, 11).Value = -0.04166000770222517

$ws.Cells.Item(5, 2).Value = 0.7472976860263175
$ws.Cells.Item(5, 3).Value = 0.2931833097224661
$ws.Cells.Item(5, 4).Value = 0.4850066528825197
$ws.Cells.Item(5, 5).Value = 0.6537292374023527
$ws.Cells.Item(5, 6).Value = 0.1338542619389697
$ws.Cells.Item(5, 7).Value = 0.4197108360949934
$ws.Cells.Item(5, 8).Value = 0.251248383465597
$ws.Cells.Item(5, 9).Value = 0.553701034065628
$ws.Cells.Item(5, 10).Value = -0.05813389123226798
$ws.Cells.Item(5, 11).Value = 0.6310880987550094

$ws.Cells.Item(6, 2).Value = 0.6333810503403763
$ws.Cells.Item(6, 3).Value = 0.5603777146664568
$ws.Cells.Item(6, 4).Value = 0.4626359510893536
$ws.Cells.Item(6, 5).Value = 0.1585481745257029
$ws.Cells.Item(6, 6).Value = 0.4290902396512893
$ws.Cells.Item(6, 7).Value = 0.196437584854206
$ws.Cells.Item(6, 8).Value = 0.5348069431769582
$ws.Cells.Item(6, 9).Value = -0.07299789204731372
$ws.Cells.Item(6, 10).Value = 0.6023763881256388
$ws.Cells.Item(6, 11).Value = 0.3304428394308724

$ws.Cells.Item(7, 2).Value = 1.011072586656707
$ws.Cells.Item(7, 3).Value = 0.509746653301712
$ws.Cells.Item(7, 4).Value = -0.08173505626510533
$ws.Cells.Item(7, 5).Value = 0.4650353207171521
$ws.Cells.Item(7, 6).Value = 0.1943868274146131
$ws.Cells.Item(7, 7).Value = 0.458238703740928
$ws.Cells.Item(7, 8).Value = -0.1003972864225971
$ws.Cells.Item(7, 9).Value = 0.5760873379279133
$ws.Cells.Item(7, 10).Value = 0.2870455487175398
$ws.Cells.Item(7, 11).ClearContents()

$ws.Cells.Item(8, 2).Value = 0.8220722402252505
$ws.Cells.Item(8, 3).Value = 0.051553683470419
$ws.Cells.Item(8, 4).Value = 0.2859568989009605
$ws.Cells.Item(8, 5).Value = 0.2225810661851803
$ws.Cells.Item(8, 6).Value = 0.4942640149333215
$ws.Cells.Item(8, 7).Value = -0.1379388566268107
$ws.Cells.Item(8, 8).Value = 0.5705449988405521
$ws.Cells.Item(8, 9).Value = 0.2912208776562884
$ws.Cells.Item(8, 10).ClearContents()

$ws.Cells.Item(9, 2).Value = 0.2871441745782602
$ws.Cells.Item(9, 3).Value = 0.3706356397752701
$ws.Cells.Item(9, 4).Value = 0.0769806995514632
$ws.Cells.Item(9, 5).Value = 0.5042177577925642
$ws.Cells.Item(9, 6).Value = -0.1225736869272658
$ws.Cells.Item(9, 7).Value = 0.5316966638831291
$ws.Cells.Item(9, 8).Value = 0.2753750686291025
$ws.Cells.Item(9, 9).ClearContents()

$ws.Cells.Item(10, 2).Value = 0.6816598262566529
$ws.Cells.Item(10, 3).Value = 0.1940694739626584
$ws.Cells.Item(10, 4).Value = 0.3417536163429973
$ws.Cells.Item(10, 5).Value = -0.0936147364620642
$ws.Cells.Item(10, 6).Value = 0.5674395363380327
$ws.Cells.Item(10, 7).Value = 0.2440474222454754
$ws.Cells.Item(10, 8).ClearContents()

$ws.Cells.Item(11, 2).Value = 0.4408946513667728
$ws.Cells.Item(11, 3).Value = 0.3594094838808868
$ws.Cells.Item(11, 4).Value = -0.1882369755730587
$ws.Cells.Item(11, 5).Value = 0.5996569245865127
$ws.Cells.Item(11, 6).Value = 0.2564355480731927
$ws.Cells.Item(11, 7).ClearContents()

$ws.Cells.Item(12, 2).Value = 0.5990858432970987
$ws.Cells.Item(12, 3).Value = -0.1032127321038452
$ws.Cells.Item(12, 4).Value = 0.4831723462284986
$ws.Cells.Item(12, 5).Value = 0.2715408197250452
$ws.Cells.Item(12, 6).ClearContents()

$ws.Cells.Item(13, 2).Value = 0.0616473449302421
$ws.Cells.Item(13, 3).Value = 0.4967096184764148
$ws.Cells.Item(13, 4).Value = 0.2085679007350822
$ws.Cells.Item(13, 5).ClearContents()

$ws.Cells.Item(14, 2).Value = 0.7505586603418228
$ws.Cells.Item(14, 3).Value = 0.3078859509171186
$ws.Cells.Item(14, 4).ClearContents()

$ws.Cells.Item(15, 2).Value = 0.3519456421565676
$ws.Cells.Item(15, 3).ClearContents()

$ws.Cells.Item(16, 2).ClearContents()
